# ArbeitsProtokoll update — "Updated Arbeitsprotokoll, new Controller"
#
# 1) Row 6: tag the "Datenbank-Verbindung" entry with a new "Teil" label
#    (K6 = "IMP…Implementierung").
# 2) Rows 18/19: effort ("Zeit an AP") revised from 2 -> 3 hours.
# 3) Rows 20-33: 14 new work-package entries appended (TP + IMP block).
# 4) The old "Weiteres:" rows (25/26) are pushed down to rows 46/47 to make
#    room for the new block, content unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: a date-formatted cell (style index used by existing column G
# cells) to stamp onto freshly written G-column dates. We copy the
# number format from an existing date cell and paste-special (formats
# only) so we reuse the workbook's existing date style instead of
# minting a new one.
# ---------------------------------------------------------------------
function Set-DateStyle($rangeAddress) {
    $ws.Range("G6").Copy() | Out-Null
    $ws.Range($rangeAddress).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 1) New K6 cell - label for row 6 ("IMP…Implementierung")
# ---------------------------------------------------------------------
$ws.Range("K6").Value = "IMP…Implementierung"

# ---------------------------------------------------------------------
# 2) Rows 18 & 19 - effort corrected from 2 to 3 hours
# ---------------------------------------------------------------------
$ws.Range("E18").Value = 3
$ws.Range("E19").Value = 3

# ---------------------------------------------------------------------
# 3) Preserve the "Weiteres:" rows (currently 25/26) before the new
#    block shoves them down to 46/47.
# ---------------------------------------------------------------------
$ws.Range("A25").Value = "Weiteres:"
$ws.Range("D25").Value = "Gespräch mit Auftraggeber"
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 43385
Set-DateStyle("G25")

$ws.Range("D26").Value = "Gespräch mit Auftraggeber und Betreuer"
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 43426
Set-DateStyle("G26")

$ws.Range("A46").Value = $ws.Range("A25").Value
$ws.Range("D46").Value = $ws.Range("D25").Value
$ws.Range("E46").Value = $ws.Range("E25").Value
$ws.Range("F46").Value = $ws.Range("F25").Value
$ws.Range("G46").Value = $ws.Range("G25").Value
Set-DateStyle("G46")

$ws.Range("D47").Value = $ws.Range("D26").Value
$ws.Range("E47").Value = $ws.Range("E26").Value
$ws.Range("F47").Value = $ws.Range("F26").Value
$ws.Range("G47").Value = $ws.Range("G26").Value
Set-DateStyle("G47")

# now clear the old 25/26 locations (their data now lives at 46/47)
$ws.Rows("25:26").Delete() | Out-Null

# ---------------------------------------------------------------------
# 4) New work-package rows 20-33
# ---------------------------------------------------------------------

# Row 20 - TP #23 Tests
$ws.Range("B20").Value = "TP"
$ws.Range("C20").Value = 23
$ws.Range("D20").Value = "Tests"
$ws.Range("E20").Value = 1
$ws.Range("G20").Value = 43436
Set-DateStyle("G20")

# Row 21 - TP #24 Server/Client-Beschreibung
$ws.Range("B21").Value = "TP"
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = "Server/Client-Beschreibung"
$ws.Range("E21").Value = 2
$ws.Range("G21").Value = 43436
Set-DateStyle("G21")

# Row 22 - TP #25 Projektstrukturplan (PSP)
$ws.Range("B22").Value = "TP"
$ws.Range("C22").Value = 25
$ws.Range("D22").Value = "Projektstrukturplan (PSP)"
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 43437
Set-DateStyle("G22")

# Row 23 - IMP Datenbank-Erstellung
$ws.Range("B23").Value = "IMP"
$ws.Range("D23").Value = "Datenbank-Erstellung"
$ws.Range("E23").Value = 4
$ws.Range("G23").Value = 43444
Set-DateStyle("G23")

# Row 24 - IMP Datenbank-Verbindung/Deployment
$ws.Range("B24").Value = "IMP"
$ws.Range("D24").Value = "Datenbank-Verbindung/Deployment"
$ws.Range("E24").Value = 4
$ws.Range("G24").Value = 43456
Set-DateStyle("G24")

# Row 25 - IMP View-Login
$ws.Range("B25").Value = "IMP"
$ws.Range("D25").Value = "View-Login"
$ws.Range("F25").Value = 3
$ws.Range("G25").Value = 43457
Set-DateStyle("G25")

# Row 26 - IMP View-Register
$ws.Range("B26").Value = "IMP"
$ws.Range("D26").Value = "View-Register"
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 43457
Set-DateStyle("G26")

# Row 27 - IMP Model-Datenbank-Operationen (laufend erweitert)
$ws.Range("B27").Value = "IMP"
$ws.Range("D27").Value = "Model-Datenbank-Operationen (laufend erweitert)"
$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 43458
Set-DateStyle("G27")

# Row 28 - IMP Model-DBConnector
$ws.Range("B28").Value = "IMP"
$ws.Range("D28").Value = "Model-DBConnector"
$ws.Range("E28").Value = 2
$ws.Range("G28").Value = 43458
Set-DateStyle("G28")

# Row 29 - IMP Controller-Login
$ws.Range("B29").Value = "IMP"
$ws.Range("D29").Value = "Controller-Login"
$ws.Range("E29").Value = 2
$ws.Range("G29").Value = 43459
Set-DateStyle("G29")

# Row 30 - IMP Controller-Register
$ws.Range("B30").Value = "IMP"
$ws.Range("D30").Value = "Controller-Register"
$ws.Range("E30").Value = 2
$ws.Range("G30").Value = 43460
Set-DateStyle("G30")

# Row 31 - IMP Controller-Resets (4 zusammengefasst)
$ws.Range("B31").Value = "IMP"
$ws.Range("D31").Value = "Controller-Resets (4 zusammengefasst)"
$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 43460
Set-DateStyle("G31")

# Row 32 - IMP Controller-MailSender
$ws.Range("B32").Value = "IMP"
$ws.Range("D32").Value = "Controller-MailSender"
$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 43461
Set-DateStyle("G32")

# Row 33 - IMP Controller-Logout
$ws.Range("B33").Value = "IMP"
$ws.Range("D33").Value = "Controller-Logout"
$ws.Range("E33").Value = 1
$ws.Range("G33").Value = 43461
Set-DateStyle("G33")

# ---------------------------------------------------------------------
# 5) Selection / cursor position like in the authored file
# ---------------------------------------------------------------------
$ws.Range("F28").Select()
